$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cellD = $ws.Range("D2")
$cellD.NumberFormat = "@"
$cellD.Value = "39.975.14"
$cellD.ClearFormats()
$ws.Range("E2").Value = "  +1.89%  "

$cellD = $ws.Range("D3")
$cellD.NumberFormat = "@"
$cellD.Value = "2.237.56"
$cellD.ClearFormats()
$ws.Range("E3").Value = "  -0.14%  "

$ws.Range("E4").Value = "  +0.09%  "

$cellD = $ws.Range("D5")
$cellD.NumberFormat = "@"
$cellD.Value = "292.54"
$cellD.ClearFormats()
$ws.Range("E5").Value = "  -1.90%  "

$cellD = $ws.Range("D6")
$cellD.NumberFormat = "@"
$cellD.Value = "86.35"
$cellD.ClearFormats()
$ws.Range("E6").Value = "  +6.85%  "

$ws.Range("E7").Value = "  +0.96%  "

$ws.Range("E8").Value = "  +0.05%  "

$cellD = $ws.Range("D9")
$cellD.NumberFormat = "@"
$cellD.Value = "0.471"
$cellD.ClearFormats()
$ws.Range("E9").Value = "  +2.57%  "

$cellD = $ws.Range("D10")
$cellD.NumberFormat = "@"
$cellD.Value = "30.93"
$cellD.ClearFormats()
$ws.Range("E10").Value = "  +10.45%  "

$cellD = $ws.Range("D11")
$cellD.NumberFormat = "@"
$cellD.Value = "0.0799"
$cellD.ClearFormats()
$ws.Range("E11").Value = "  +2.93%  "

$cellD = $ws.Range("D12")
$cellD.NumberFormat = "@"
$cellD.Value = "47.07"
$cellD.ClearFormats()
$ws.Range("E12").Value = "  +1.14%  "

$ws.Range("E13").Value = "  +0.24%  "

$cellD = $ws.Range("D14")
$cellD.NumberFormat = "@"
$cellD.Value = "6.42"
$cellD.ClearFormats()
$ws.Range("E14").Value = "  +5.46%  "

$cellD = $ws.Range("D15")
$cellD.NumberFormat = "@"
$cellD.Value = "2.584.40"
$cellD.ClearFormats()
$ws.Range("E15").Value = "  +0.16%  "

$cellD = $ws.Range("D16")
$cellD.NumberFormat = "@"
$cellD.Value = "14.19"
$cellD.ClearFormats()
$ws.Range("E16").Value = "  +0.85%  "

$cellD = $ws.Range("D17")
$cellD.NumberFormat = "@"
$cellD.Value = "2.242.04"
$cellD.ClearFormats()
$ws.Range("E17").Value = "  +0.00%  "

$ws.Range("E18").Value = "  +2.14%  "

$cellD = $ws.Range("D19")
$cellD.NumberFormat = "@"
$cellD.Value = "39.901.89"
$cellD.ClearFormats()
$ws.Range("E19").Value = "  +2.05%  "

$cellD = $ws.Range("D20")
$cellD.NumberFormat = "@"
$cellD.Value = "0.0₃0893"
$cellD.ClearFormats()
$ws.Range("E20").Value = "  +3.80%  "

$ws.Range("E21").Value = "  +1.01%  "

$cellD = $ws.Range("D22")
$cellD.NumberFormat = "@"
$cellD.Value = "65.65"
$cellD.ClearFormats()
$ws.Range("E22").Value = "  +0.54%  "

$cellD = $ws.Range("D23")
$cellD.NumberFormat = "@"
$cellD.Value = "10.53"
$cellD.ClearFormats()
$ws.Range("E23").Value = "  +5.93%  "

$cellD = $ws.Range("D24")
$cellD.NumberFormat = "@"
$cellD.Value = "235.82"
$cellD.ClearFormats()
$ws.Range("E24").Value = "  +3.97%  "

$ws.Range("E25").Value = "  +0.03%  "

$cellD = $ws.Range("D26")
$cellD.NumberFormat = "@"
$cellD.Value = "2.44"
$cellD.ClearFormats()
$ws.Range("E26").Value = "  +1.85%  "

$ws.Range("E27").Value = "  +5.86%  "

$cellD = $ws.Range("D28")
$cellD.NumberFormat = "@"
$cellD.Value = "23.06"
$cellD.ClearFormats()
$ws.Range("E28").Value = "  +3.27%  "

$cellD = $ws.Range("D29")
$cellD.NumberFormat = "@"
$cellD.Value = "2.23"
$cellD.ClearFormats()
$ws.Range("E29").Value = "  +2.23%  "

$ws.Range("E30").Value = "  +4.23%  "

$cellD = $ws.Range("D31")
$cellD.NumberFormat = "@"
$cellD.Value = "34.00"
$cellD.ClearFormats()
$ws.Range("E31").Value = "  +7.28%  "

$cellD = $ws.Range("D32")
$cellD.NumberFormat = "@"
$cellD.Value = "153.82"
$cellD.ClearFormats()
$ws.Range("E32").Value = "  +3.82%  "

$ws.Range("E33").Value = "  +0.14%  "

$ws.Range("E34").Value = "  +2.19%  "

$ws.Range("E35").Value = "  +4.18%  "

$ws.Range("E36").Value = "  +2.27%  "

$cellD = $ws.Range("D37")
$cellD.NumberFormat = "@"
$cellD.Value = "16.45"
$cellD.ClearFormats()
$ws.Range("E37").Value = "  +11.28%  "

$ws.Range("E38").Value = "  +1.83%  "

$ws.Range("E39").Value = "  +2.69%  "

$ws.Range("E40").Value = "  +2.77%  "

$ws.Range("E41").Value = "  +3.89%  "

$ws.Range("E42").Value = "  +3.89%  "

$cellD = $ws.Range("D43")
$cellD.NumberFormat = "@"
$cellD.Value = "1.973.30"
$cellD.ClearFormats()
$ws.Range("E43").Value = "  +2.99%  "

$cellD = $ws.Range("D44")
$cellD.NumberFormat = "@"
$cellD.Value = "2.23"
$cellD.ClearFormats()
$ws.Range("E44").Value = "  +2.62%  "

$ws.Range("E45").Value = "  +6.90%  "

$cellD = $ws.Range("D46")
$cellD.NumberFormat = "@"
$cellD.Value = "9.81"
$cellD.ClearFormats()
$ws.Range("E46").Value = "  +9.43%  "

$cellD = $ws.Range("D47")
$cellD.NumberFormat = "@"
$cellD.Value = "16.23"
$cellD.ClearFormats()
$ws.Range("E47").Value = "  -0.95%  "

$ws.Range("E48").Value = "  +1.67%  "

$cellD = $ws.Range("D49")
$cellD.NumberFormat = "@"
$cellD.Value = "2.455.06"
$cellD.ClearFormats()
$ws.Range("E49").Value = "  +0.12%  "

$ws.Range("E50").Value = "  +6.64%  "

$ws.Range("E51").Value = "  +13.59%  "

